$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.919.32'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '2.836.85'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '358.77'
$ws.Range('E5').Value = '  +4.05%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '113.05'
$ws.Range('E6').Value = '  -2.13%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.568'
$ws.Range('E7').Value = '  +3.87%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.601'
$ws.Range('E9').Value = '  +4.01%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.37'
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0861'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '20.17'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.82'
$ws.Range('E14').Value = '  +2.62%  '
$ws.Range('D15').Value = '3.291.70'
$ws.Range('E15').Value = '  +2.34%  '
$ws.Range('D16').Value = '2.832.88'
$ws.Range('E16').Value = '  +1.96%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.928'
$ws.Range('E17').Value = '  +5.79%  '
$ws.Range('D18').Value = '51.859.43'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.57'
$ws.Range('E19').Value = '  +7.85%  '
$ws.Range('E20').Value = '  -1.12%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.49'
$ws.Range('E21').Value = '  +1.83%  '
$ws.Range('D22').Value = '0.0₃0993'
$ws.Range('E22').Value = '  +1.61%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.22'
$ws.Range('E23').Value = '  +0.35%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '269.54'
$ws.Range('E24').Value = '  -2.48%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.84'
$ws.Range('E25').Value = '  +2.40%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '27.09'
$ws.Range('E26').Value = '  +1.51%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.36'
$ws.Range('E28').Value = '  +1.88%  '
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '53.25'
$ws.Range('E30').Value = '  +5.89%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.140'
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '35.47'
$ws.Range('E32').Value = '  +2.22%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0473'
$ws.Range('E33').Value = '  +23.29%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.93'
$ws.Range('E34').Value = '  +3.71%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.47'
$ws.Range('E35').Value = '  +10.63%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0858'
$ws.Range('E36').Value = '  +4.76%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.30'
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.06'
$ws.Range('E39').Value = '  -1.84%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.64'
$ws.Range('E40').Value = '  -1.77%  '
$ws.Range('E41').Value = '  +1.19%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '23.56'
$ws.Range('E42').Value = '  +2.15%  '
$ws.Range('E43').Value = '  -4.47%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '125.04'
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('E45').Value = '  -3.47%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.39'
$ws.Range('E46').Value = '  +2.03%  '
$ws.Range('D47').Value = '2.106.84'
$ws.Range('E47').Value = '  +1.77%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.99'
$ws.Range('E49').Value = '  +7.87%  '
$ws.Range('E50').Value = '  +11.47%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '61.78'
$ws.Range('E51').Value = '  +3.99%  '
